$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "GRT-USD"
